## VehicleLabExperimentFiles.xlsx — "update with guinardian 5ml vert runs"
##
## 1. Rename the second sheet ("Sheet1") to "9um_beads".
## 2. Add a new sheet "Guinardia" after it, with the 5ml vertical-run data
##    for the Guinardia experiment (5 rows + header).
## 3. Make "Guinardia" the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- 1. rename existing "Sheet1" -> "9um_beads" --------------------------
$beads = $wb.Worksheets.Item(2)
$beads.Name = "9um_beads"

# Row 1 on the beads sheet had the header row selected (A1:XFD1) rather
# than the previous cell-selection - mirrors the "select whole header row"
# step noted in the source diff.
$beads.Activate()
$beads.Rows.Item(1).Select()

# --- 2. add the new "Guinardia" sheet right after "9um_beads" ------------
$gui = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $beads)
$gui.Name = "Guinardia"

# Headers (same columns/shared strings as the "9um_beads" sheet)
$gui.Range("A1").Value = "Filename"
$gui.Range("B1").Value = "Volume"
$gui.Range("C1").Value = "HorzOrVert"
$gui.Range("D1").Value = "CellConc"
$gui.Range("E1").Value = "ml_analyzed"
$gui.Range("F1").Value = "runtime"
$gui.Range("G1").Value = "inhibittime"
$gui.Range("H1").Value = "numtriggers"
$gui.Range("I1").Value = "Comments"
$gui.Range("J1").Value = "Comments2"

# header styles match "9um_beads": A1 uses the highlighted-fill style,
# D1:H1 use the centered-number style
$gui.Range("A1").Style = $beads.Range("A1").Style
$gui.Range("D1:H1").Style = $beads.Range("D1").Style

# Data rows
$data = @(
    @("D20151112T164628", 5, "V", 373.436, 4.1212, 1198, 193.9, 1539, "use all triggers", "2nd syringe from gui, still some beads, first sucked 5ml sample then run ~1ml"),
    @("D20151112T170931", 5, "V", 492,     3.9977, 1198, 224.2, 1967, "use all triggers", "still some beads"),
    @("D20151112T173151", 5, "V", 629.58,  3.8645, 1198, 262.5, 2433, "use all triggers", "still some beads"),
    @("D20151112T175410", 5, "V", 742.7,   3.7457, 1198, 284,   2782, "use all triggers", "still some beads"),
    @("D20151112T181728", 5, "V", 790.7,   3.7144, 1198, 291.2, 2937, "use all triggers", "still some beads")
)

$r = 2
foreach ($row in $data) {
    $gui.Cells.Item($r, 1).Value = $row[0]
    $gui.Cells.Item($r, 2).Value = $row[1]
    $gui.Cells.Item($r, 3).Value = $row[2]
    $gui.Cells.Item($r, 4).Value = $row[3]
    $gui.Cells.Item($r, 5).Value = $row[4]
    $gui.Cells.Item($r, 6).Value = $row[5]
    $gui.Cells.Item($r, 7).Value = $row[6]
    $gui.Cells.Item($r, 8).Value = $row[7]
    $gui.Cells.Item($r, 9).Value = $row[8]
    $gui.Cells.Item($r, 10).Value = $row[9]
    $r = $r + 1
}

# column widths to match the "9um_beads" sheet layout
$gui.Columns.Item(1).ColumnWidth = 18.42578125
$gui.Columns.Item(3).ColumnWidth = 11
$gui.Columns.Item(4).ColumnWidth = 8.7109375
$gui.Columns.Item(5).ColumnWidth = 12.140625
$gui.Columns.Item(7).ColumnWidth = 11
$gui.Columns.Item(8).ColumnWidth = 11.7109375

# freeze the header row, scroll/select like the source file (I6:J6 selected
# on the last data row) and make this the active/selected sheet
$gui.Activate()
$gui.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$gui.Range("I6:J6").Select()
